# Refined metadata to be additional tab
#
# 1) Update the "time_taken" (column F) timestamps on the existing "data"
#    sheet for every data row (2..112) to the new query run's timestamps.
# 2) Add a new "metadata" worksheet (placed after "data") describing the
#    panel query itself (name/id/version/version_created/query_time/request).

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

$newTimestamps = @(
    "2021-10-05 14:20:04.799086",
    "2021-10-05 14:20:04.799094",
    "2021-10-05 14:20:04.799098",
    "2021-10-05 14:20:04.799100",
    "2021-10-05 14:20:04.799103",
    "2021-10-05 14:20:04.799106",
    "2021-10-05 14:20:04.799108",
    "2021-10-05 14:20:04.799111",
    "2021-10-05 14:20:04.799113",
    "2021-10-05 14:20:04.799116",
    "2021-10-05 14:20:04.799119",
    "2021-10-05 14:20:04.799121",
    "2021-10-05 14:20:04.799124",
    "2021-10-05 14:20:04.799126",
    "2021-10-05 14:20:04.799129",
    "2021-10-05 14:20:04.799131",
    "2021-10-05 14:20:04.799134",
    "2021-10-05 14:20:04.799137",
    "2021-10-05 14:20:04.799139",
    "2021-10-05 14:20:04.799142",
    "2021-10-05 14:20:04.799144",
    "2021-10-05 14:20:04.799147",
    "2021-10-05 14:20:04.799149",
    "2021-10-05 14:20:04.799152",
    "2021-10-05 14:20:04.799154",
    "2021-10-05 14:20:04.799157",
    "2021-10-05 14:20:04.799160",
    "2021-10-05 14:20:04.799162",
    "2021-10-05 14:20:04.799165",
    "2021-10-05 14:20:04.799167",
    "2021-10-05 14:20:04.799170",
    "2021-10-05 14:20:04.799172",
    "2021-10-05 14:20:04.799175",
    "2021-10-05 14:20:04.799178",
    "2021-10-05 14:20:04.799180",
    "2021-10-05 14:20:04.799183",
    "2021-10-05 14:20:04.799185",
    "2021-10-05 14:20:04.799188",
    "2021-10-05 14:20:04.799190",
    "2021-10-05 14:20:04.799193",
    "2021-10-05 14:20:04.799195",
    "2021-10-05 14:20:04.799198",
    "2021-10-05 14:20:04.799201",
    "2021-10-05 14:20:04.799203",
    "2021-10-05 14:20:04.799206",
    "2021-10-05 14:20:04.799208",
    "2021-10-05 14:20:04.799211",
    "2021-10-05 14:20:04.799213",
    "2021-10-05 14:20:04.799216",
    "2021-10-05 14:20:04.799218",
    "2021-10-05 14:20:04.799221",
    "2021-10-05 14:20:04.799223",
    "2021-10-05 14:20:04.799226",
    "2021-10-05 14:20:04.799229",
    "2021-10-05 14:20:04.799231",
    "2021-10-05 14:20:04.799234",
    "2021-10-05 14:20:04.799236",
    "2021-10-05 14:20:04.799239",
    "2021-10-05 14:20:04.799241",
    "2021-10-05 14:20:04.799244",
    "2021-10-05 14:20:04.799246",
    "2021-10-05 14:20:04.799249",
    "2021-10-05 14:20:04.799251",
    "2021-10-05 14:20:04.799253",
    "2021-10-05 14:20:04.799257",
    "2021-10-05 14:20:04.799260",
    "2021-10-05 14:20:04.799263",
    "2021-10-05 14:20:04.799265",
    "2021-10-05 14:20:04.799268",
    "2021-10-05 14:20:04.799270",
    "2021-10-05 14:20:04.799273",
    "2021-10-05 14:20:04.799275",
    "2021-10-05 14:20:04.799277",
    "2021-10-05 14:20:04.799280",
    "2021-10-05 14:20:04.799283",
    "2021-10-05 14:20:04.799285",
    "2021-10-05 14:20:04.799289",
    "2021-10-05 14:20:04.799292",
    "2021-10-05 14:20:04.799295",
    "2021-10-05 14:20:04.799297",
    "2021-10-05 14:20:04.799300",
    "2021-10-05 14:20:04.799302",
    "2021-10-05 14:20:04.799305",
    "2021-10-05 14:20:04.799307",
    "2021-10-05 14:20:04.799310",
    "2021-10-05 14:20:04.799312",
    "2021-10-05 14:20:04.799315",
    "2021-10-05 14:20:04.799317",
    "2021-10-05 14:20:04.799320",
    "2021-10-05 14:20:04.799322",
    "2021-10-05 14:20:04.799325",
    "2021-10-05 14:20:04.799327",
    "2021-10-05 14:20:04.799331",
    "2021-10-05 14:20:04.799334",
    "2021-10-05 14:20:04.799337",
    "2021-10-05 14:20:04.799339",
    "2021-10-05 14:20:04.799342",
    "2021-10-05 14:20:04.799344",
    "2021-10-05 14:20:04.799347",
    "2021-10-05 14:20:04.799349",
    "2021-10-05 14:20:04.799352",
    "2021-10-05 14:20:04.799354",
    "2021-10-05 14:20:04.799357",
    "2021-10-05 14:20:04.799359",
    "2021-10-05 14:20:04.799362",
    "2021-10-05 14:20:04.799364",
    "2021-10-05 14:20:04.799367",
    "2021-10-05 14:20:04.799369",
    "2021-10-05 14:20:04.799374",
    "2021-10-05 14:20:04.799377",
    "2021-10-05 14:20:04.799380"
)

$startRow = 2
for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $startRow + $i
    $dataSheet.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# --- Add the "metadata" sheet, positioned right after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Reuse the bold/bordered/centered header style from data!B1 for the header
# row (B1:G1) and from data!A2 for the index cell (A2), matching the look
# of the "data" sheet.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)

$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Early onset dystonia"
$metaSheet.Range("C2").Value = 192
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.89"
$metaSheet.Range("E2").Value = "2021-08-20T13:48:46.341314Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:20:04.795733"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/192/?format=json"

# Keep "data" as the active sheet/tab, matching the original activeTab=0.
$dataSheet.Activate()
